$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.814.14"
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.953.42"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.87"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.84"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.952.41"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +2.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +5.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  +5.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.02"
$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.444.23"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.693.57"
$ws.Range("E17").Value = "  +2.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.975.86"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "441.62"
$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.668"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.07"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.43"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("E26").Value = "  -3.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.85"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  +4.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.61"
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0970"
$ws.Range("E32").Value = "  +10.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.48"
$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.50"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.51"
$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -4.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.93"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.733.70"
$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.49"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "365.29"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0338"
$ws.Range("E48").Value = "  -2.74%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.85"
$ws.Range("E51").Value = "  -3.93%  "
